$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$ws.Range("N15").Value = "月極め、春/冬/春休み短期プログラム"

# Row 18
$ws.Range("G18").Value = "香川県立中央病院　院長　高口　浩一"
$ws.Range("O18").Value = "構造：木造1階建て（専用建物）、保育室：27.75㎡、屋外：133.73㎡135.90㎡屋外：あり"

# Row 19
$ws.Range("O19").Value = "構造：鉄筋コンクリート造4階建ての2階（病院建物）、保育室：188.6㎡、屋外：264㎡"

# Row 20
$ws.Range("G20").Value = "高松赤十字病院　院長　西村　健夫"

# Row 23
$ws.Range("I23").Value = "ー"
$ws.Range("J23").Value = "ー"
$ws.Range("K23").Value = "ー"
$ws.Range("L23").Value = "ー"

# Row 24
$ws.Range("J24").Value = "8：00～18：00利用児がいるときのみ開所"

# Row 32
$ws.Range("H32").Value = "R4.6.1"
$ws.Range("I32").Value = "8：45～16：45"
$ws.Range("J32").Value = "休"
$ws.Range("K32").Value = "休"
$ws.Range("L32").Value = "'5"

# Row 33 - L33 numeric-looking text, keep as text
$ws.Range("L33").Value = "'6"

# Row 34
$ws.Range("F34").Value = "087-802-5360"
$ws.Range("H34").Value = "R4.4.1"
$ws.Range("I34").Value = "7：30～78：30（18：30～19：00）"
$ws.Range("J34").Value = "7：30～18：30（18：30～19：00）"
$ws.Range("K34").Value = "休"
$ws.Range("L34").Value = "'60"
$ws.Range("M34").Value = "3～5歳児"
$ws.Range("N34").Value = "月極め"
$ws.Range("O34").Value = "構造：鉄骨造2階建ての1階（専用建物）、保育室：132.84㎡、屋外：なし"
